$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.763.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.663.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.89%  "

$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3650"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3237"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.140"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07053"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.066"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.661.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.602"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001049"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06610"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.99%  "

$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.929"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.776.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.453"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.406"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.845.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.77%  "

$ws.Range("E31").Value = "  +2.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.069"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.688"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08483"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.637"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.149"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02258"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06028"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.220"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2079"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.202"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("E43").Value = "  +0.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5915"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.846"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5651"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.946"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06960"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.185"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.77%  "
